# The deck currently carries the "Integral" (Red Violet) colour theme on
# its slide master/design (ppt/theme/theme2.xml), while the unused
# ppt/theme/theme1.xml still holds the stock "Office Theme" palette.
# The target revision swaps the two so the design in use reverts to the
# default Office colour palette (dk1/lt1/dk2/lt2/accent1-6/hlink/folHlink).
#
# PowerPoint's Design/Theme object model doesn't give script access to
# raw part names, so we apply the visible, scriptable effect of that
# swap: push the stock Office theme colours into the presentation's
# live ThemeColorScheme (MsoThemeColorSchemeIndex 1-12, in
# dk1,lt1,dk2,lt2,accent1..accent6,hlink,folHlink order).

$p = $ppt.ActivePresentation

function RGBVal($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0,2),16)
    $g = [Convert]::ToInt32($hex.Substring(2,2),16)
    $b = [Convert]::ToInt32($hex.Substring(4,2),16)
    return $r + ($g * 256) + ($b * 65536)
}

# Office Theme colour scheme, in MsoThemeColorSchemeIndex order.
$officeColors = @(
    "000000",  # 1  dk1
    "FFFFFF",  # 2  lt1
    "44546A",  # 3  dk2
    "E7E6E6",  # 4  lt2
    "5B9BD5",  # 5  accent1
    "ED7D31",  # 6  accent2
    "A5A5A5",  # 7  accent3
    "FFC000",  # 8  accent4
    "4472C4",  # 9  accent5
    "70AD47",  # 10 accent6
    "0563C1",  # 11 hlink
    "954F72"   # 12 folHlink
)

$tcs = $p.Slides.Item(1).ThemeColorScheme
for ($i = 1; $i -le 12; $i++) {
    $tcs.Item($i).RGB = RGBVal($officeColors[$i - 1])
}
